$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Existing-cell edits ---
$ws.Range("O1140").Value = 2
$ws.Range("R1142").Value = 0
$ws.Range("R1143").Value = 0

# --- Append 13 new weekly rows (1144-1156) ---
# Row 1144
$ws.Range("A1144").Value = 45474
$ws.Range("B1144").Value = 3884
$ws.Range("C1144").Value = 4047.35009765625
$ws.Range("D1144").Value = 3884
$ws.Range("E1144").Value = 4011.800048828125
$ws.Range("F1144").Value = 4002.50390625
$ws.Range("G1144").Value = 10973987
$ws.Range("H1144").Value = 2024
$ws.Range("I1144").Value = 7
$ws.Range("J1144").Value = 1
$ws.Range("K1144").Value = 0
$ws.Range("L1144").Value = 0
$ws.Range("M1144").Value = 0
$ws.Range("N1144").Value = 27
$ws.Range("O1144").Value = 0
$ws.Range("P1144").Value = 0
$ws.Range("Q1144").Value = 0
$ws.Range("A1144").NumberFormat = "YYYY-MM-DD HH:MM:SS"

# Row 1145
$ws.Range("A1145").Value = 45481
$ws.Range("B1145").Value = 4022
$ws.Range("C1145").Value = 4199.9501953125
$ws.Range("D1145").Value = 3895.60009765625
$ws.Range("E1145").Value = 4183.9501953125
$ws.Range("F1145").Value = 4174.25537109375
$ws.Range("G1145").Value = 24115752
$ws.Range("H1145").Value = 2024
$ws.Range("I1145").Value = 7
$ws.Range("J1145").Value = 8
$ws.Range("K1145").Value = 0
$ws.Range("L1145").Value = 0
$ws.Range("M1145").Value = 0
$ws.Range("N1145").Value = 28
$ws.Range("O1145").Value = 0
$ws.Range("P1145").Value = 0
$ws.Range("Q1145").Value = 0
$ws.Range("A1145").NumberFormat = "YYYY-MM-DD HH:MM:SS"

# Row 1146
$ws.Range("A1146").Value = 45488
$ws.Range("B1146").Value = 4235.7001953125
$ws.Range("C1146").Value = 4358.75
$ws.Range("D1146").Value = 4144.89990234375
$ws.Range("E1146").Value = 4302.39990234375
$ws.Range("F1146").Value = 4292.4306640625
$ws.Range("G1146").Value = 16394269
$ws.Range("H1146").Value = 2024
$ws.Range("I1146").Value = 7
$ws.Range("J1146").Value = 15
$ws.Range("K1146").Value = 0
$ws.Range("L1146").Value = 0
$ws.Range("M1146").Value = 0
$ws.Range("N1146").Value = 29
$ws.Range("O1146").Value = 0
$ws.Range("P1146").Value = 0
$ws.Range("Q1146").Value = 1
$ws.Range("A1146").NumberFormat = "YYYY-MM-DD HH:MM:SS"

# Row 1147
$ws.Range("A1147").Value = 45495
$ws.Range("B1147").Value = 4299.9501953125
$ws.Range("C1147").Value = 4422.4501953125
$ws.Range("D1147").Value = 4265
$ws.Range("E1147").Value = 4387.85009765625
$ws.Range("F1147").Value = 4387.85009765625
$ws.Range("G1147").Value = 12433939
$ws.Range("H1147").Value = 2024
$ws.Range("I1147").Value = 7
$ws.Range("J1147").Value = 22
$ws.Range("K1147").Value = 0
$ws.Range("L1147").Value = 0
$ws.Range("M1147").Value = 0
$ws.Range("N1147").Value = 30
$ws.Range("O1147").Value = 0
$ws.Range("P1147").Value = 0
$ws.Range("Q1147").Value = 0
$ws.Range("A1147").NumberFormat = "YYYY-MM-DD HH:MM:SS"

# Row 1148
$ws.Range("A1148").Value = 45502
$ws.Range("B1148").Value = 4415
$ws.Range("C1148").Value = 4431
$ws.Range("D1148").Value = 4270
$ws.Range("E1148").Value = 4283.0498046875
$ws.Range("F1148").Value = 4283.0498046875
$ws.Range("G1148").Value = 10888326
$ws.Range("H1148").Value = 2024
$ws.Range("I1148").Value = 7
$ws.Range("J1148").Value = 29
$ws.Range("K1148").Value = 0
$ws.Range("L1148").Value = 0
$ws.Range("M1148").Value = 0
$ws.Range("N1148").Value = 31
$ws.Range("O1148").Value = 0
$ws.Range("P1148").Value = 0
$ws.Range("Q1148").Value = 0
$ws.Range("A1148").NumberFormat = "YYYY-MM-DD HH:MM:SS"

# Row 1149
$ws.Range("A1149").Value = 45509
$ws.Range("B1149").Value = 4239
$ws.Range("C1149").Value = 4258.25
$ws.Range("D1149").Value = 4110.5
$ws.Range("E1149").Value = 4228.75
$ws.Range("F1149").Value = 4228.75
$ws.Range("G1149").Value = 9972323
$ws.Range("H1149").Value = 2024
$ws.Range("I1149").Value = 8
$ws.Range("J1149").Value = 5
$ws.Range("K1149").Value = 0
$ws.Range("L1149").Value = 0
$ws.Range("M1149").Value = 0
$ws.Range("N1149").Value = 32
$ws.Range("O1149").Value = 0
$ws.Range("P1149").Value = 0
$ws.Range("Q1149").Value = 0
$ws.Range("A1149").NumberFormat = "YYYY-MM-DD HH:MM:SS"

# Row 1150
$ws.Range("A1150").Value = 45516
$ws.Range("B1150").Value = 4230
$ws.Range("C1150").Value = 4427
$ws.Range("D1150").Value = 4183
$ws.Range("E1150").Value = 4416.0498046875
$ws.Range("F1150").Value = 4416.0498046875
$ws.Range("G1150").Value = 9116108
$ws.Range("H1150").Value = 2024
$ws.Range("I1150").Value = 8
$ws.Range("J1150").Value = 12
$ws.Range("K1150").Value = 0
$ws.Range("L1150").Value = 0
$ws.Range("M1150").Value = 0
$ws.Range("N1150").Value = 33
$ws.Range("O1150").Value = 0
$ws.Range("P1150").Value = 0
$ws.Range("Q1150").Value = 0
$ws.Range("A1150").NumberFormat = "YYYY-MM-DD HH:MM:SS"

# Row 1151
$ws.Range("A1151").Value = 45523
$ws.Range("B1151").Value = 4421
$ws.Range("C1151").Value = 4565
$ws.Range("D1151").Value = 4390.35009765625
$ws.Range("E1151").Value = 4463.89990234375
$ws.Range("F1151").Value = 4463.89990234375
$ws.Range("G1151").Value = 9803623
$ws.Range("H1151").Value = 2024
$ws.Range("I1151").Value = 8
$ws.Range("J1151").Value = 19
$ws.Range("K1151").Value = 0
$ws.Range("L1151").Value = 0
$ws.Range("M1151").Value = 0
$ws.Range("N1151").Value = 34
$ws.Range("O1151").Value = 0
$ws.Range("P1151").Value = 0
$ws.Range("Q1151").Value = 0
$ws.Range("A1151").NumberFormat = "YYYY-MM-DD HH:MM:SS"

# Row 1152
$ws.Range("A1152").Value = 45530
$ws.Range("B1152").Value = 4489
$ws.Range("C1152").Value = 4592.25
$ws.Range("D1152").Value = 4457.2998046875
$ws.Range("E1152").Value = 4553.75
$ws.Range("F1152").Value = 4553.75
$ws.Range("G1152").Value = 10393782
$ws.Range("H1152").Value = 2024
$ws.Range("I1152").Value = 8
$ws.Range("J1152").Value = 26
$ws.Range("K1152").Value = 0
$ws.Range("L1152").Value = 0
$ws.Range("M1152").Value = 0
$ws.Range("N1152").Value = 35
$ws.Range("O1152").Value = 0
$ws.Range("P1152").Value = 0
$ws.Range("Q1152").Value = 0
$ws.Range("A1152").NumberFormat = "YYYY-MM-DD HH:MM:SS"

# Row 1153
$ws.Range("A1153").Value = 45537
$ws.Range("B1153").Value = 4576
$ws.Range("C1153").Value = 4588
$ws.Range("D1153").Value = 4436.7001953125
$ws.Range("E1153").Value = 4456.75
$ws.Range("F1153").Value = 4456.75
$ws.Range("G1153").Value = 7537240
$ws.Range("H1153").Value = 2024
$ws.Range("I1153").Value = 9
$ws.Range("J1153").Value = 2
$ws.Range("K1153").Value = 0
$ws.Range("L1153").Value = 0
$ws.Range("M1153").Value = 0
$ws.Range("N1153").Value = 36
$ws.Range("O1153").Value = 0
$ws.Range("P1153").Value = 0
$ws.Range("Q1153").Value = 0
$ws.Range("A1153").NumberFormat = "YYYY-MM-DD HH:MM:SS"

# Row 1154
$ws.Range("A1154").Value = 45544
$ws.Range("B1154").Value = 4455
$ws.Range("C1154").Value = 4549.35009765625
$ws.Range("D1154").Value = 4430.5
$ws.Range("E1154").Value = 4522.60009765625
$ws.Range("F1154").Value = 4522.60009765625
$ws.Range("G1154").Value = 8097450
$ws.Range("H1154").Value = 2024
$ws.Range("I1154").Value = 9
$ws.Range("J1154").Value = 9
$ws.Range("K1154").Value = 0
$ws.Range("L1154").Value = 0
$ws.Range("M1154").Value = 0
$ws.Range("N1154").Value = 37
$ws.Range("O1154").Value = 0
$ws.Range("P1154").Value = 0
$ws.Range("Q1154").Value = 0
$ws.Range("A1154").NumberFormat = "YYYY-MM-DD HH:MM:SS"

# Row 1155
$ws.Range("A1155").Value = 45551
$ws.Range("B1155").Value = 4527
$ws.Range("C1155").Value = 4546.0498046875
$ws.Range("D1155").Value = 4227.5
$ws.Range("E1155").Value = 4284.89990234375
$ws.Range("F1155").Value = 4284.89990234375
$ws.Range("G1155").Value = 14191722
$ws.Range("H1155").Value = 2024
$ws.Range("I1155").Value = 9
$ws.Range("J1155").Value = 16
$ws.Range("K1155").Value = 0
$ws.Range("L1155").Value = 0
$ws.Range("M1155").Value = 0
$ws.Range("N1155").Value = 38
$ws.Range("O1155").Value = 0
$ws.Range("P1155").Value = 0
$ws.Range("Q1155").Value = 0
$ws.Range("A1155").NumberFormat = "YYYY-MM-DD HH:MM:SS"

# Row 1156
$ws.Range("A1156").Value = 45558
$ws.Range("B1156").Value = 4300
$ws.Range("C1156").Value = 4378
$ws.Range("D1156").Value = 4240
$ws.Range("E1156").Value = 4308.7001953125
$ws.Range("F1156").Value = 4308.7001953125
$ws.Range("G1156").Value = 14404146
$ws.Range("H1156").Value = 2024
$ws.Range("I1156").Value = 9
$ws.Range("J1156").Value = 23
$ws.Range("K1156").Value = 0
$ws.Range("L1156").Value = 0
$ws.Range("M1156").Value = 0
$ws.Range("N1156").Value = 39
$ws.Range("O1156").Value = 0
$ws.Range("P1156").Value = 0
$ws.Range("Q1156").Value = 0
$ws.Range("A1156").NumberFormat = "YYYY-MM-DD HH:MM:SS"
